$d = $word.ActiveDocument

# Original paragraph text is "hjj" in a single run.
# Target: three runs -> "H", "jj", "n,n, ;nbc c" (same fr-FR formatting).

# 1) Capitalize the initial "h" -> "H" (still within the first run for now).
$r1 = $d.Range(0, 1)
$r1.Text = "H"

# 2) Append the new trailing text after "hjj" (now "Hjj").
$r3 = $d.Range(3, 3)
$r3.InsertAfter("n,n, ;nbc c")

# Touch formatting on the appended text so Word keeps it as its own run
# instead of silently re-merging it with the preceding run.
$r3after = $d.Range(3, 14)
$r3after.Font.Bold = $true
$r3after.Font.Bold = $false

# 3) Touch formatting on the leading "H" the same way so it stays split
# from the following "jj" run.
$r1b = $d.Range(0, 1)
$r1b.Font.Bold = $true
$r1b.Font.Bold = $false
